# Auto update Excel log
# Append 7 new PRESENCE_DETECTED rows to the "mmWave" sheet (rows 46-52)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

$newRows = @(
    @("2026-02-01", "17:41:14", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:41:20", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:41:31", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:41:41", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:41:52", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:42:02", "17:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "17:42:12", "17:00", "Living Room", "PRESENCE_DETECTED", "Active")
)

$startRow = 46
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    # Column A looks like an ISO date ("2026-02-01"); Excel would normally
    # auto-convert such text into a date serial number on assignment. Force
    # it to be stored as plain text (matching the rest of the log) by
    # stamping the cell as Text first, then restore the default "Normal"
    # style so no stray number-format style lingers on the cell.
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $rowData[0]
    $dateCell.Style = "Normal"

    for ($c = 2; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}
